$wb = $excel.ActiveWorkbook

# --- 1. Amount sheet: bump the balance from 10000 to 10001 ---
$wsAmount = $wb.Worksheets.Item("Amount")
$wsAmount.Range("A1").Value = 10001

# --- 2. Withdraw History: append the (blank) transaction row that was being
#        skipped, then land the selection on the new row (fixes the "writes
#        one cell ahead" bug: cursor now sits on the just-written row, not
#        past it) ---
$wsWithdraw = $wb.Worksheets.Item("Withdraw History")
$wsWithdraw.Range("A2").Value = "tmp"
$wsWithdraw.Range("A2").ClearContents()
$wsWithdraw.Range("A2").Style = "Normal"
$wsWithdraw.Range("A2:D2").Select()

# --- 3. Deposit History: same fix ---
$wsDeposit = $wb.Worksheets.Item("Deposit History")
$wsDeposit.Range("A2").Value = "tmp"
$wsDeposit.Range("A2").ClearContents()
$wsDeposit.Range("A2").Style = "Normal"
$wsDeposit.Range("A2:E2").Select()

# --- 4. Transfer History: just land the cursor where the next entry would go ---
$wsTransfer = $wb.Worksheets.Item("Transfer History")
$wsTransfer.Range("A13").Select()

# --- 5. Absolute History: land the cursor here too, and leave this as the
#        active sheet/tab when the workbook is saved ---
$wsAbsolute = $wb.Worksheets.Item("Absolute History")
$wsAbsolute.Activate()
$wsAbsolute.Range("A9").Select()
